$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (item id 5484)
$ws.Range("H32").Value = 905.625
$ws.Range("I32").Value = 749.8333
$ws.Range("J32").Value = 999.1
$ws.Range("K32").Value = 749.8333
$ws.Range("L32").Value = 999.1
$ws.Range("M32").Value = -423.8333
$ws.Range("N32").Value = -1651.1

# Row 33 (item id 5512)
$ws.Range("H33").Value = 441.33334
$ws.Range("I33").Value = 258.44446
$ws.Range("K33").Value = 258.44446
$ws.Range("M33").Value = -29.44445999999999

# Row 40 (item id 5505)
$ws.Range("H40").Value = 3501
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825

# Row 70 (item id 12604)
$ws.Range("H70").Value = 5041.1665
$ws.Range("J70").Value = 5041.1665
$ws.Range("L70").Value = 15123.4995
$ws.Range("N70").Value = -15663.4995

# Row 73 (item id 12604)
$ws.Range("H73").Value = 5041.1665
$ws.Range("J73").Value = 5041.1665
$ws.Range("L73").Value = 15123.4995
$ws.Range("N73").Value = -16995.4995

# Row 100 (item id 19906)
$ws.Range("H100").Value = 3012
$ws.Range("I100").Value = 675
$ws.Range("J100").Value = 5349
$ws.Range("K100").Value = 675
$ws.Range("L100").Value = 5349
$ws.Range("M100").Value = -134
$ws.Range("N100").Value = -6431

# Row 111 (item id 27768)
$ws.Range("H111").Value = 1860
$ws.Range("I111").Value = 895
$ws.Range("J111").Value = 2825
$ws.Range("K111").Value = 2685
$ws.Range("L111").Value = 8475
$ws.Range("M111").Value = 382
$ws.Range("N111").Value = -14609

# Row 112 (item id 27960)
$ws.Range("H112").Value = 2602.9565
$ws.Range("J112").Value = 2689.0476
$ws.Range("L112").Value = 8067.1428
$ws.Range("N112").Value = -10283.1428

# Row 113 (item id 27775)
$ws.Range("H113").Value = 3976.375
$ws.Range("I113").Value = 3398.8
$ws.Range("J113").Value = 4939
$ws.Range("K113").Value = 3398.8
$ws.Range("L113").Value = 4939
$ws.Range("M113").Value = -144.8000000000002
$ws.Range("N113").Value = -11447

# Row 129 (item id 36115)
$ws.Range("H129").Value = 2995.2144
$ws.Range("I129").Value = 1210.8
$ws.Range("J129").Value = 3986.5557
$ws.Range("K129").Value = 3632.4
$ws.Range("L129").Value = 11959.6671
$ws.Range("M129").Value = 1367.6
$ws.Range("N129").Value = -21959.6671

# Row 137 (item id 44013)
$ws.Range("H137").Value = 1805.9166
$ws.Range("I137").Value = 1863
$ws.Range("J137").Value = 1691.75
$ws.Range("K137").Value = 5589
$ws.Range("L137").Value = 5075.25
$ws.Range("M137").Value = -3039
$ws.Range("N137").Value = -10175.25

# Row 141 (item id 44161)
$ws.Range("H141").Value = 3550.375
$ws.Range("I141").Value = 3953.1428
$ws.Range("K141").Value = 11859.4284
$ws.Range("M141").Value = -6679.428400000001

$ws = $wb.Worksheets.Item("ARM")
# Row 63 (item id 12528)
$ws.Range("H63").Value = 2409.75
$ws.Range("J63").Value = 3600
$ws.Range("L63").Value = 3600
$ws.Range("N63").Value = -4972

# Row 66 (item id 12528)
$ws.Range("H66").Value = 2409.75
$ws.Range("J66").Value = 3600
$ws.Range("L66").Value = 18000
$ws.Range("N66").Value = -24864

# Row 74 (item id 44000)
$ws.Range("H74").Value = 1345.4615
$ws.Range("I74").Value = 1332.5834
$ws.Range("K74").Value = 1332.5834
$ws.Range("M74").Value = -458.5834

# Row 77 (item id 44000)
$ws.Range("H77").Value = 1345.4615
$ws.Range("I77").Value = 1332.5834
$ws.Range("K77").Value = 6662.916999999999
$ws.Range("M77").Value = -2294.916999999999

# Row 132 (item id 43997)
$ws.Range("H132").Value = 1136.8667
$ws.Range("I132").Value = 1136.8667
$ws.Range("K132").Value = 3410.6001
$ws.Range("M132").Value = -880.6001000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 64 (item id 14184)
$ws.Range("H64").Value = 35715104
$ws.Range("J64").Value = 55556028
$ws.Range("L64").Value = 55556028
$ws.Range("N64").Value = -55556478

# Row 67 (item id 14184)
$ws.Range("H67").Value = 35715104
$ws.Range("J67").Value = 55556028
$ws.Range("L67").Value = 55556028
$ws.Range("N67").Value = -55557588

# Row 105 (item id 19947)
$ws.Range("H105").Value = 2200.5
$ws.Range("J105").Value = 2226.5
$ws.Range("L105").Value = 2226.5
$ws.Range("N105").Value = -5720.5

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (item id 5367)
$ws.Range("H22").Value = 188
$ws.Range("I22").Value = 174.57143
$ws.Range("J22").Value = 250.66667
$ws.Range("K22").Value = 174.57143
$ws.Range("L22").Value = 250.66667
$ws.Range("M22").Value = 175.42857
$ws.Range("N22").Value = -950.6666700000001

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (item id 4650)
$ws.Range("H4").Value = 2059.524
$ws.Range("I4").Value = 1960
$ws.Range("J4").Value = 2308.3333
$ws.Range("K4").Value = 5880
$ws.Range("L4").Value = 6924.999899999999
$ws.Range("M4").Value = -5768
$ws.Range("N4").Value = -7148.999899999999

# Row 7 (item id 4728)
$ws.Range("H7").Value = 357
$ws.Range("I7").Value = 447.5
$ws.Range("J7").Value = 320.8
$ws.Range("K7").Value = 1342.5
$ws.Range("L7").Value = 962.4000000000001
$ws.Range("M7").Value = -1230.5
$ws.Range("N7").Value = -1186.4

# Row 17 (item id 4640)
$ws.Range("H17").Value = 15024.5
$ws.Range("I17").Value = 50
$ws.Range("K17").Value = 150
$ws.Range("M17").Value = 19

# Row 46 (item id 4701)
$ws.Range("H46").Value = 2997.5
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15182

# Row 59 (item id 4694)
$ws.Range("H59").Value = 1905
$ws.Range("I59").Value = 1905
$ws.Range("K59").Value = 5715
$ws.Range("M59").Value = -5175

# Row 61 (item id 4727)
$ws.Range("H61").Value = 401.8889
$ws.Range("I61").Value = 349.5
$ws.Range("J61").Value = 506.66666
$ws.Range("K61").Value = 1048.5
$ws.Range("L61").Value = 1519.99998
$ws.Range("M61").Value = -833.5
$ws.Range("N61").Value = -1949.99998

# Row 70 (item id 12867)
$ws.Range("H70").Value = 2777
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73 (item id 12867)
$ws.Range("H73").Value = 2777
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (item id 44008)
$ws.Range("H132").Value = 1599.6666
$ws.Range("I132").Value = 1599.6666
$ws.Range("K132").Value = 4798.9998
$ws.Range("M132").Value = -2268.9998

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (item id 5282)
$ws.Range("H46").Value = 2956.2856
$ws.Range("I46").Value = 2699.375
$ws.Range("J46").Value = 3778.4
$ws.Range("K46").Value = 2699.375
$ws.Range("L46").Value = 3778.4
$ws.Range("M46").Value = -2511.375
$ws.Range("N46").Value = -4154.4

# Row 100 (item id 19995)
$ws.Range("H100").Value = 1499
$ws.Range("I100").Value = 1998.5
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 1998.5
$ws.Range("L100").Value = 500
$ws.Range("M100").Value = -1457.5
$ws.Range("N100").Value = -1582

# Row 132 (item id 44058)
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (item id 44029)
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
